# TC02_Canine_Filter_SamplePatho-BCellLymphoma.xlsx
# Commit: "Fixed Tests for SamplePatholoy, SampleType, Se, StageOfDisease, and Study"
#
# The stored Cypher query for the "CasesTab" row (column B, row 2 of the
# "startup" sheet) used to OPTIONAL MATCH a (:cohort) node and return an
# extra `Cohort` column. That column was dropped, so the query text is
# updated to match: the trailing
#     coalesce(co.cohort_description, '') AS `Cohort`
# line is removed, and the comma that used to separate it from the
# preceding `Response to Treatment` line is removed too.
#
# (The SamplesTab / FilesTab query text in B3 / B4 is unchanged -- only the
# shared-string table slot backing them shifts because B2's text changed.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = 'MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (samp:sample)-->(c) 
 WHERE samp.specific_sample_pathology IN ["B Cell Lymphoma"]  
OPTIONAL MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '''') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '''') AS `Study Code` ,
        coalesce(s.clinical_study_type, '''') AS  `Study Type`,
        coalesce(demo.breed, '''') AS Breed ,
        coalesce(diag.disease_term, '''') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '''') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '''') AS Age ,
        coalesce(demo.sex, '''') AS Sex ,
        coalesce(demo.neutered_indicator, '''') AS `Neutered Status`,
        coalesce(demo.weight, '''') AS `Weight (kg)`,
        coalesce(diag.best_response, '''') AS `Response to Treatment`'

$ws.Range("B2").Value = $casesQuery

# Keep the saved selection in sync with the workbook (B2, not the old B4).
$ws.Range("B2").Select() | Out-Null
